$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"-0.001038695952956914"
$ws.Range("D2").Value = [double]"-0.01600057219184237"
$ws.Range("E2").Value = [double]"-0.01383385043300223"
$ws.Range("F2").Value = [double]"-0.2977330752182752"
$ws.Range("G2").Value = [double]"0.002916539466241375"
$ws.Range("H2").Value = [double]"-0.05027020654233638"
$ws.Range("I2").Value = [double]"0.0009328710948466323"
$ws.Range("C3").Value = [double]"-0.1902780891161342"
$ws.Range("D3").Value = [double]"-3.651395549482004"
$ws.Range("E3").Value = [double]"-2.191468958199948"
$ws.Range("F3").Value = [double]"-39.23138677816314"
$ws.Range("H3").Value = [double]"2.728091192841021"
$ws.Range("C4").Value = [double]"-0.003455267146364349"
$ws.Range("D4").Value = [double]"-0.07892941749014426"
$ws.Range("E4").Value = [double]"-0.04767234976316104"
$ws.Range("F4").Value = [double]"-1.118406115856487"
$ws.Range("G4").Value = [double]"2.213805419160053e-05"
$ws.Range("H4").Value = [double]"-0.1750327392073814"
$ws.Range("I4").Value = [double]"-8.300338231492788e-05"
$ws.Range("J4").Value = [double]"-0.1557759204224567"
$ws.Range("D5").Value = [double]"-0.01065126268525773"
$ws.Range("E5").Value = [double]"-0.01427961863228688"
$ws.Range("F5").Value = [double]"-0.2253112770381449"
$ws.Range("H5").Value = [double]"-0.03339100210155266"
$ws.Range("J5").Value = [double]"-0.02382061940124913"
$ws.Range("C6").Value = [double]"-0.002259173036691209"
$ws.Range("D6").Value = [double]"-0.02546117937163217"
$ws.Range("E6").Value = [double]"-0.02383432661008555"
$ws.Range("I6").Value = [double]"-0.07400143065024167"
$ws.Range("J6").Value = [double]"-0.005386661128795822"
$ws.Range("C7").Value = [double]"3.042886328330496e-05"
$ws.Range("D7").Value = [double]"2.220135915464994e-05"
$ws.Range("E7").Value = [double]"0.0001120885785610426"
$ws.Range("I7").Value = [double]"0.02070887600712012"
$ws.Range("J7").Value = [double]"0.002185896000355569"
$ws.Range("C8").Value = [double]"-0.007532573839171164"
$ws.Range("D8").Value = [double]"-0.005481143451561366"
$ws.Range("E8").Value = [double]"-0.02762085277481674"
$ws.Range("I8").Value = [double]"-0.07367226153291995"
$ws.Range("J8").Value = [double]"-0.006921796556753179"
$ws.Range("C9").Value = [double]"-0.007329303108889462"
$ws.Range("D9").Value = [double]"-0.005347078649890591"
$ws.Range("E9").Value = [double]"-0.02698654027113889"
$ws.Range("I9").Value = [double]"-0.1685449741066805"
$ws.Range("C10").Value = [double]"-0.002606588845992519"
$ws.Range("D10").Value = [double]"-0.001510081983724376"
$ws.Range("E10").Value = [double]"-0.0003108644982603437"
$ws.Range("I10").Value = [double]"-0.01134794281824725"
$ws.Range("J10").Value = [double]"-0.0009170765320050123"
$ws.Range("C11").Value = [double]"-0.03355990815907717"
$ws.Range("D11").Value = [double]"0.8868028017459437"
$ws.Range("E11").Value = [double]"-0.1207279756199569"
$ws.Range("I11").Value = [double]"-1.128603667253628"
$ws.Range("J11").Value = [double]"-0.02519565066904761"
